# This script swaps a set of field values between row 2 and row 4 of the
# "Artfynd" sheet, matching the data correction described in the diff
# (two observation records had gotten their Id/coordinates/time/biotope-
# description/substrate fields crossed, so those columns are exchanged
# back between the two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "Q", "R", "Z", "AB", "AI", "AJ", "AK", "AO")

foreach ($col in $columns) {
    $cell2 = $ws.Range("$col`2")
    $cell4 = $ws.Range("$col`4")

    $val2 = $cell2.Value2
    $val4 = $cell4.Value2

    $cell2.Value2 = $val4
    $cell4.Value2 = $val2
}
